$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '37.519.63'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +5.34%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.057.79'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +3.79%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '255.06'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +3.88%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.657'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +3.43%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '66.55'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +11.81%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.398'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +9.46%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '60.60'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +2.59%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0793'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +7.14%  '
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +0.63%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.941'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -1.12%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '23.86'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +27.57%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.17'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +3.54%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.360.49'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +3.88%  '
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +8.21%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.057.38'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +3.41%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '37.462.04'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +5.28%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '73.88'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +3.07%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0₃0895'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +5.31%  '
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +5.92%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '242.26'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +3.80%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.65'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +3.85%  '
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.04%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +7.62%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.10'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +9.50%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '162.81'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -1.51%  '
$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '20.27'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +5.10%  '
$ws.Range('B30').Value = 'Kaspa'
$ws.Range('C30').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.135'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +40.48%  '
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +3.58%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.23'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +7.01%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.21'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +8.14%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0636'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +6.33%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +8.14%  '
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +0.35%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.38'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +15.33%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +0.17%  '
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +2.53%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.12'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +35.83%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.104'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +12.97%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.29'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +4.75%  '
$ws.Range('B43').Value = 'HuobiToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.07'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +7.72%  '
$ws.Range('B44').Value = 'InjectiveProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '18.04'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +9.59%  '
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +6.83%  '
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +3.49%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '97.31'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +3.96%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.04'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +2.77%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.419.84'
$ws.Range('B50').Value = 'MultiversX'
$ws.Range('C50').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '49.05'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +4.22%  '
$ws.Range('B51').Value = 'MXToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.95'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +2.19%  '
